$d = $word.ActiveDocument

# Avoid Word's "smart quotes" auto-substitution from mangling apostrophes
# during Find/Replace and Range.Text assignment.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ---------------------------------------------------------------------
# Edit 1: "Pour pouvoir emballer des articles , votre transfert ..."
#         -> "Pour pouvoir emballer des articles, votre transfert ..."
#         (drop the stray space before the comma, and the proofErr wrap)
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("des articles , votre")
$r.Text = "des articles, votre"

# Re-split "articles," back into its own run (matches the original
# run layout minus the w:proofErr markers).
$r1 = $d.Content
[void]$r1.Find.Execute("articles,")
$r1.Font.Size = 12.5
$r1.Font.Size = 12

# ---------------------------------------------------------------------
# Edit 2: "Dans l'onglet Opérations,ouvez mettre vos articles ..."
#         -> "Dans l'onglet Opérations, vous pouvez mettre vos articles ..."
#         (fix the typo/missing "vous p" and drop the proofErr wraps)
# ---------------------------------------------------------------------
$r2 = $d.Content
[void]$r2.Find.Execute("Dans l'onglet Opérations,ouvez mettre")
$r2.Text = "Dans l'onglet Opérations, vous pouvez mettre"

# Re-split the merged run into "Dans l'onglet Opérations," / " vous pouvez"
# (the trailing " mettre vos articles..." run is left untouched).
$r3 = $d.Content
[void]$r3.Find.Execute(" vous pouvez")
$r3.Font.Size = 12.5
$r3.Font.Size = 12

# ---------------------------------------------------------------------
# Edit 3: "Répétez les mêmes étapes pour les autres article que vous ..."
#         -> "... pour les autres articles que vous ..."
#         (add the missing "s" and drop the proofErr wrap)
# ---------------------------------------------------------------------
$r4 = $d.Content
[void]$r4.Find.Execute("pour les autres article que")
$r4.Text = "pour les autres articles que"

# Re-split "les autres articles" back into its own run.
$r5 = $d.Content
[void]$r5.Find.Execute("les autres articles")
$r5.Font.Size = 12.5
$r5.Font.Size = 12
